# Generate Report for Handoff
#
# Moves the localization status from "In Translation" to "Ready for
# handoff" and refreshes the "Latest HO Xliff Generate Date" / "Latest
# Handoff Datetime" timestamps, then widens the status-datetime columns
# slightly to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed timestamps ----------------------------------------------
$overview.Range("G2").Value = "2016-08-25 04:38:10"
$zhcn.Range("H2").Value     = "2016-08-25 04:38:03"
$dede.Range("H2").Value     = "2016-08-25 04:38:10"

# --- Widen the status/datetime columns to fit the new text -----------------
# (the engine stores column width quantised to 1/6-character pixel steps,
# so 16.38 is the input that lands on the closest representable width)
$overview.Range("E1:F1").ColumnWidth = 16.38
$zhcn.Range("C1").ColumnWidth = 16.38
$dede.Range("C1").ColumnWidth = 16.38
